$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.05045533333333333
$ws.Range("H2").Value = 0.151366
$ws.Range("I2").Value = 0.004442474524580737
$ws.Range("J2").Value = 0.004442474524580737
$ws.Range("M2").Value = 4.093680666666667
$ws.Range("N2").Value = 12.281042
$ws.Range("O2").Value = 0.1610908176055751
$ws.Range("P2").Value = 0.161090817605575
$ws.Range("Q2").Value = 0.2065480225968889
$ws.Range("R2").Value = 1.858932203372
$ws.Range("S2").Value = 0.0007156418533566494
$ws.Range("T2").Value = 0.0007156418533566493
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.05045533333333333
$ws.Range("H3").Value = 0.151366
$ws.Range("I3").Value = 0.004442474524580737
$ws.Range("J3").Value = 0.004442474524580737
$ws.Range("O3").Value = 0.5606512265211691
$ws.Range("P3").Value = 0.5606512265211691
$ws.Range("Q3").Value = 0.7188578711419999
$ws.Range("R3").Value = 6.469720840278
$ws.Range("S3").Value = 0.002490678790995238
$ws.Range("T3").Value = 0.002490678790995238
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.05045533333333333
$ws.Range("H4").Value = 0.151366
$ws.Range("I4").Value = 0.004442474524580737
$ws.Range("J4").Value = 0.004442474524580737
$ws.Range("M4").Value = 7.071161666666666
$ws.Range("N4").Value = 21.213485
$ws.Range("O4").Value = 0.2782579558732559
$ws.Range("P4").Value = 0.2782579558732559
$ws.Range("Q4").Value = 0.3567778189455555
$ws.Range("R4").Value = 3.21100037051
$ws.Range("S4").Value = 0.00123615388022885
$ws.Range("T4").Value = 0.00123615388022885
$ws.Range("I5").Value = 0.7425623198471305
$ws.Range("J5").Value = 0.7425623198471305
$ws.Range("M5").Value = 4.093680666666667
$ws.Range("N5").Value = 12.281042
$ws.Range("O5").Value = 0.1610908176055751
$ws.Range("P5").Value = 0.161090817605575
$ws.Range("Q5").Value = 34.52462765306645
$ws.Range("R5").Value = 310.721648877598
$ws.Range("S5").Value = 0.1196199712272668
$ws.Range("T5").Value = 0.1196199712272668
$ws.Range("I6").Value = 0.7425623198471305
$ws.Range("J6").Value = 0.7425623198471305
$ws.Range("O6").Value = 0.5606512265211691
$ws.Range("P6").Value = 0.5606512265211691
$ws.Range("S6").Value = 0.4163184753906983
$ws.Range("T6").Value = 0.4163184753906983
$ws.Range("I7").Value = 0.7425623198471305
$ws.Range("J7").Value = 0.7425623198471305
$ws.Range("M7").Value = 7.071161666666666
$ws.Range("N7").Value = 21.213485
$ws.Range("O7").Value = 0.2782579558732559
$ws.Range("P7").Value = 0.2782579558732559
$ws.Range("Q7").Value = 59.63562952141277
$ws.Range("R7").Value = 536.720665692715
$ws.Range("S7").Value = 0.2066238732291653
$ws.Range("T7").Value = 0.2066238732291653
$ws.Range("G8").Value = 2.873389
$ws.Range("H8").Value = 8.620167
$ws.Range("I8").Value = 0.2529952056282888
$ws.Range("J8").Value = 0.2529952056282888
$ws.Range("M8").Value = 4.093680666666667
$ws.Range("N8").Value = 12.281042
$ws.Range("O8").Value = 0.1610908176055751
$ws.Range("P8").Value = 0.161090817605575
$ws.Range("Q8").Value = 11.76273699711267
$ws.Range("R8").Value = 105.864632974014
$ws.Range("S8").Value = 0.04075520452495163
$ws.Range("T8").Value = 0.04075520452495162
$ws.Range("G9").Value = 2.873389
$ws.Range("H9").Value = 8.620167
$ws.Range("I9").Value = 0.2529952056282888
$ws.Range("J9").Value = 0.2529952056282888
$ws.Range("O9").Value = 0.5606512265211691
$ws.Range("P9").Value = 0.5606512265211691
$ws.Range("Q9").Value = 40.938354045879
$ws.Range("R9").Value = 368.445186412911
$ws.Range("S9").Value = 0.1418420723394755
$ws.Range("T9").Value = 0.1418420723394755
$ws.Range("G10").Value = 2.873389
$ws.Range("H10").Value = 8.620167
$ws.Range("I10").Value = 0.2529952056282888
$ws.Range("J10").Value = 0.2529952056282888
$ws.Range("M10").Value = 7.071161666666666
$ws.Range("N10").Value = 21.213485
$ws.Range("O10").Value = 0.2782579558732559
$ws.Range("P10").Value = 0.2782579558732559
$ws.Range("Q10").Value = 20.31819815022167
$ws.Range("R10").Value = 182.863783351995
$ws.Range("S10").Value = 0.07039792876386168
$ws.Range("T10").Value = 0.07039792876386168
